# Add explicit root to TemplateExcelFileGenerator
#
# - Row 2 (A2) text gains a trailing "which is:" clause.
# - A brand-new row is inserted right after it (becomes row 3) holding the
#   absolute root-folder path, merged across A3:F3 just like rows 1-2.
# - The pre-existing "Pol0_90" / "Pol45_135" row is pushed down to row 5
#   (a blank row 4 separates it from the new header rows, same as before).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the explanatory text in row 2.
$ws.Range("A2").Value = "The path must be the remaining path after the root folder, which is:"

# 2. Insert a fresh row at position 3; this pushes the old row 3 (blank) and
#    row 4 (Pol0_90 / Pol45_135) down by one, landing the data row on row 5.
$ws.Rows.Item(3).Insert()

# 3. Populate & merge the new row 3 with the explicit root folder path.
$ws.Range("A3:F3").Merge()
$ws.Range("A3").Value = "/home/masoud/Documents/four-polar/fourPolar-io/target/test-classes/fr/fresnel/fourPolar/io/imageSet/acquisition/sample/finders/excel"
$ws.Range("B3:F3").ClearFormats()

Write-Output "Applied explicit-root edit to Sheet0"
